$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 26 ("RM 232") entirely - the remaining rows shift up.
$ws.Rows.Item(26).Delete()

# After the shift above, the old row 28 ("SC 92") is now row 27;
# remove it too so the rest of the rows shift up again.
$ws.Rows.Item(27).Delete()

# Re-run the missing-data imputation pass: some previously-blank cells
# now carry a computed value, and some previously-filled cells are
# cleared back to blank (column D on the untouched rows, plus columns
# B and D on the rows that shifted up from the block below).
$ws.Range("D6").Value = -14.2
$ws.Range("D8").ClearContents()
$ws.Range("D19").Value = -15.5
$ws.Range("D21").ClearContents()
$ws.Range("D23").Value = -13.9

$ws.Range("B26").ClearContents()
$ws.Range("B27").Value = -20.4
$ws.Range("D27").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("D29").Value = -13
